$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the changed values in rows 2-5 (columns C and E)
$ws.Range("C2").Value = 23.26790258615237
$ws.Range("E2").Value = 0.5691875071735751

$ws.Range("C3").Value = 23.95886060057096
$ws.Range("E3").Value = 0.6050178600059178

$ws.Range("C4").Value = 25.33142633381867
$ws.Range("E4").Value = 0.677828129208542

$ws.Range("C5").Value = 28.71055627543781
$ws.Range("E5").Value = 0.823414760904444

# Remove rows 6 through 17 so the used range shrinks to A1:E5
$ws.Rows("6:17").Delete()
